# sql_lib.xlsx: save rpg and sql and clone
# - Rename two "Read Active ... Account" rows to the 定期/活期 variants
# - Append a new row (35) with a "BCM staff" entry and its SQL text
# - Update the selected/active cell to B35

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename existing category labels (row 16 = checking/活期, row 17 = time deposit/定期)
$ws.Range("B16").Value = "Read Active 活期 Account w/ Personal Information"
$ws.Range("B17").Value = "Read Active 定期 Account w/ Personal Information"

# Add the new row 35
$ws.Range("A35").Value = "SQL"
$ws.Range("B35").Value = "BCM staff "
$ws.Range("C35").Value = "select cup00301.cunbr                 `n from zusrlib/cup00301, zusrlib/cup027                `nwhere cup00301.cunbr=cup027.cunbr and cuten1 like 'B%'"

# Match formatting used by the other data rows
$ws.Rows.Item(35).RowHeight = 37.5
$ws.Range("C35").WrapText = $true

# Update view state: scroll + selection on the new row
$ws.Activate()
$ws.Range("B35").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 2
